$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (header row 1 untouched), columns A..T.
# A/B/C/D are text labels (shared strings); E..T are numeric.
$data = @(
    @("ECs", "Sertad1", "Ar", "ECs", 3.0, 1.0, 27.70991333333333, 83.12974, 0.499916436778856, 0.4999164367788561, 2.0, 0.6666666666666666, 0.480739, 1.442217, 0.01989190335837768, 0.01989190335837768, 13.32123602595333, 119.89112423358, 0.009944289447669528, 0.00994428944766953),
    @("ECs", "Sertad1", "Ar", "FAPs", 3.0, 1.0, 27.70991333333333, 83.12974, 0.499916436778856, 0.4999164367788561, 3.0, 1.0, 19.80485166666667, 59.41455500000001, 0.8194804153196195, 0.8194804153196195, 548.7907232628555, 4939.1165093657, 0.4096717292366412, 0.4096717292366413),
    @("ECs", "Sertad1", "Ar", "sCs", 3.0, 1.0, 27.70991333333333, 83.12974, 0.499916436778856, 0.4999164367788561, 3.0, 1.0, 3.881981, 11.645943, 0.1606276813220029, 0.1606276813220029, 107.5693570716466, 968.1242136448199, 0.08030041809454529, 0.08030041809454529),
    @("FAPs", "Sertad1", "Ar", "ECs", 3.0, 1.0, 19.365057, 58.095171, 0.3493663143945636, 0.3493663143945637, 2.0, 0.6666666666666666, 0.480739, 1.442217, 0.01989190335837768, 0.01989190335837768, 9.309538137122999, 83.78584323410699, 0.006949560962609252, 0.006949560962609252),
    @("FAPs", "Sertad1", "Ar", "FAPs", 3.0, 1.0, 19.365057, 58.095171, 0.3493663143945636, 0.3493663143945637, 3.0, 1.0, 19.80485166666667, 59.41455500000001, 0.8194804153196195, 0.8194804153196195, 383.5220814015451, 3451.698732613906, 0.2862988524187418, 0.2862988524187418),
    @("FAPs", "Sertad1", "Ar", "sCs", 3.0, 1.0, 19.365057, 58.095171, 0.3493663143945636, 0.3493663143945637, 3.0, 1.0, 3.881981, 11.645943, 0.1606276813220029, 0.1606276813220029, 75.17478333791699, 676.5730500412529, 0.05611790101321264, 0.05611790101321264),
    @("sCs", "Sertad1", "Ar", "ECs", 3.0, 1.0, 8.35412, 25.06236, 0.1507172488265804, 0.1507172488265804, 2.0, 0.6666666666666666, 0.480739, 1.442217, 0.01989190335837768, 0.01989190335837768, 4.016151294679999, 36.14536165211999, 0.002998052948098898, 0.002998052948098898),
    @("sCs", "Sertad1", "Ar", "FAPs", 3.0, 1.0, 8.35412, 25.06236, 0.1507172488265804, 0.1507172488265804, 3.0, 1.0, 19.80485166666667, 59.41455500000001, 0.8194804153196195, 0.8194804153196195, 165.4521074055334, 1489.0689666498, 0.1235098336642365, 0.1235098336642365),
    @("sCs", "Sertad1", "Ar", "sCs", 3.0, 1.0, 8.35412, 25.06236, 0.1507172488265804, 0.1507172488265804, 3.0, 1.0, 3.881981, 11.645943, 0.1606276813220029, 0.1606276813220029, 32.43053511172, 291.8748160054799, 0.02420936221424496, 0.02420936221424496)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
